$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fechaValue = Get-Date -Year 2022 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0

$rows = @(
    @{ Row=155; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=$fechaValue; E=8; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Lapins';  L='Primera'; M=100; N=5000; O=6000; P=5500; Q='$/caja 10 kilos'; R='Región de Ñuble'; S=550; T=10 },
    @{ Row=156; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=$fechaValue; E=8; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Lapins';  L='Segunda'; M=50;  N=4000; O=4000; P=4000; Q='$/caja 10 kilos'; R='Región de Ñuble'; S=400; T=10 },
    @{ Row=157; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=$fechaValue; E=8; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Santina'; L='Primera'; M=100; N=5000; O=6000; P=5500; Q='$/caja 10 kilos'; R='Región de Ñuble'; S=550; T=10 },
    @{ Row=158; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=$fechaValue; E=8; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Santina'; L='Segunda'; M=50;  N=4000; O=4000; P=4000; Q='$/caja 10 kilos'; R='Región de Ñuble'; S=400; T=10 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
